# Auto-generated: updates cached market price / profit figures on several sheets
# (mirrors a scheduled data refresh; workbook has no formulas, values are static)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 680.3333
$ws.Range("J29").Value = 857.3333
$ws.Range("L29").Value = 2571.9999
$ws.Range("N29").Value = -3133.9999
$ws.Range("H38").Value = 1751.2188
$ws.Range("I38").Value = 197.66667
$ws.Range("J38").Value = 2683.35
$ws.Range("K38").Value = 593.00001
$ws.Range("L38").Value = 8050.049999999999
$ws.Range("M38").Value = -221.00001
$ws.Range("N38").Value = -8794.049999999999
$ws.Range("H39").Value = 5376610
$ws.Range("I39").Value = 122.5
$ws.Range("J39").Value = 11111530
$ws.Range("K39").Value = 367.5
$ws.Range("L39").Value = 33334590
$ws.Range("M39").Value = -71.5
$ws.Range("N39").Value = -33335182
$ws.Range("H58").Value = 2257.9092
$ws.Range("I58").Value = 232.11111
$ws.Range("J58").Value = 3660.3845
$ws.Range("K58").Value = 696.3333299999999
$ws.Range("L58").Value = 10981.1535
$ws.Range("M58").Value = -546.3333299999999
$ws.Range("N58").Value = -11281.1535
$ws.Range("H80").Value = 577.4286
$ws.Range("I80").Value = 600.25
$ws.Range("J80").Value = 568.3
$ws.Range("K80").Value = 1800.75
$ws.Range("L80").Value = 1704.9
$ws.Range("M80").Value = -802.75
$ws.Range("N80").Value = -3700.9
$ws.Range("H83").Value = 577.4286
$ws.Range("I83").Value = 600.25
$ws.Range("J83").Value = 568.3
$ws.Range("K83").Value = 5402.25
$ws.Range("L83").Value = 5114.7
$ws.Range("M83").Value = -410.25
$ws.Range("N83").Value = -15098.7
$ws.Range("H100").Value = 55556652
$ws.Range("I100").Value = 1232.5
$ws.Range("K100").Value = 1232.5
$ws.Range("M100").Value = -691.5
$ws.Range("H132").Value = 3850.2856
$ws.Range("I132").Value = 2205.3
$ws.Range("J132").Value = 13720.2
$ws.Range("K132").Value = 6615.900000000001
$ws.Range("L132").Value = 41160.60000000001
$ws.Range("M132").Value = -4085.900000000001
$ws.Range("N132").Value = -46220.60000000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 34446
$ws.Range("J44").Value = 34446
$ws.Range("L44").Value = 34446
$ws.Range("N44").Value = -35422
$ws.Range("H55").Value = 18901.5
$ws.Range("J55").Value = 18901.5
$ws.Range("L55").Value = 18901.5
$ws.Range("N55").Value = -19531.5
$ws.Range("H80").Value = 25267.4
$ws.Range("J80").Value = 25267.4
$ws.Range("L80").Value = 25267.4
$ws.Range("N80").Value = -27263.4
$ws.Range("H83").Value = 25267.4
$ws.Range("J83").Value = 25267.4
$ws.Range("L83").Value = 75802.20000000001
$ws.Range("N83").Value = -85786.20000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 27898.238
$ws.Range("J82").Value = 33342.875
$ws.Range("L82").Value = 33342.875
$ws.Range("N82").Value = -34108.875
$ws.Range("H85").Value = 27898.238
$ws.Range("J85").Value = 33342.875
$ws.Range("L85").Value = 33342.875
$ws.Range("N85").Value = -35994.875
$ws.Range("H86").Value = 2531.3333
$ws.Range("I86").Value = 2521.6667
$ws.Range("J86").Value = 2570
$ws.Range("K86").Value = 2521.6667
$ws.Range("L86").Value = 2570
$ws.Range("M86").Value = -1398.6667
$ws.Range("N86").Value = -4816
$ws.Range("H89").Value = 2531.3333
$ws.Range("I89").Value = 2521.6667
$ws.Range("J89").Value = 2570
$ws.Range("K89").Value = 12608.3335
$ws.Range("L89").Value = 12850
$ws.Range("M89").Value = -6992.333500000001
$ws.Range("N89").Value = -24082
$ws.Range("H94").Value = 507.17648
$ws.Range("I94").Value = 507.17648
$ws.Range("K94").Value = 507.17648
$ws.Range("M94").Value = -56.17648000000003
$ws.Range("H99").Value = 2298.75
$ws.Range("J99").Value = 3614.2856
$ws.Range("L99").Value = 3614.2856
$ws.Range("N99").Value = -6610.2856
$ws.Range("H107").Value = 1948.909
$ws.Range("I107").Value = 1454.2222
$ws.Range("J107").Value = 4175
$ws.Range("K107").Value = 1454.2222
$ws.Range("L107").Value = 4175
$ws.Range("M107").Value = 465.7778000000001
$ws.Range("N107").Value = -8015

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 262.5
$ws.Range("I22").Value = 250
$ws.Range("J22").Value = 300
$ws.Range("K22").Value = 250
$ws.Range("L22").Value = 300
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = -1000
$ws.Range("H31").Value = 2158.2686
$ws.Range("I31").Value = 1587.9736
$ws.Range("J31").Value = 2905.5518
$ws.Range("K31").Value = 1587.9736
$ws.Range("L31").Value = 2905.5518
$ws.Range("M31").Value = -1292.9736
$ws.Range("N31").Value = -3495.5518
$ws.Range("H34").Value = 2158.2686
$ws.Range("I34").Value = 1587.9736
$ws.Range("J34").Value = 2905.5518
$ws.Range("K34").Value = 1587.9736
$ws.Range("L34").Value = 2905.5518
$ws.Range("M34").Value = -1385.9736
$ws.Range("N34").Value = -3309.5518
$ws.Range("H41").Value = 16922.5
$ws.Range("J41").Value = 20207
$ws.Range("L41").Value = 20207
$ws.Range("N41").Value = -21063
$ws.Range("H109").Value = 10957.143
$ws.Range("J109").Value = 10957.143
$ws.Range("L109").Value = 10957.143
$ws.Range("N109").Value = -13037.143
$ws.Range("H132").Value = 1977.44
$ws.Range("I132").Value = 1953.7894
$ws.Range("K132").Value = 5861.3682
$ws.Range("M132").Value = -3331.3682

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3832.5386
$ws.Range("J43").Value = 8797.799999999999
$ws.Range("L43").Value = 8797.799999999999
$ws.Range("N43").Value = -9099.799999999999
$ws.Range("H46").Value = 9228.5
$ws.Range("J46").Value = 10022
$ws.Range("L46").Value = 10022
$ws.Range("N46").Value = -10334
$ws.Range("H57").Value = 25030.25
$ws.Range("J57").Value = 25030.25
$ws.Range("L57").Value = 25030.25
$ws.Range("N57").Value = -26670.25
$ws.Range("H70").Value = 22450.5
$ws.Range("I70").Value = 24881.451
$ws.Range("J70").Value = 4739.2856
$ws.Range("K70").Value = 24881.451
$ws.Range("L70").Value = 4739.2856
$ws.Range("M70").Value = -24611.451
$ws.Range("N70").Value = -5279.2856
$ws.Range("H73").Value = 22450.5
$ws.Range("I73").Value = 24881.451
$ws.Range("J73").Value = 4739.2856
$ws.Range("K73").Value = 24881.451
$ws.Range("L73").Value = 4739.2856
$ws.Range("M73").Value = -23945.451
$ws.Range("N73").Value = -6611.2856
$ws.Range("H97").Value = 2384.8
$ws.Range("J97").Value = 3000
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992
$ws.Range("H123").Value = 33605
$ws.Range("J123").Value = 33605
$ws.Range("L123").Value = 33605
$ws.Range("N123").Value = -38505
$ws.Range("H132").Value = 2881.6758
$ws.Range("I132").Value = 2068.45
$ws.Range("J132").Value = 3838.4119
$ws.Range("K132").Value = 6205.349999999999
$ws.Range("L132").Value = 11515.2357
$ws.Range("M132").Value = -3675.349999999999
$ws.Range("N132").Value = -16575.2357

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 18023
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 18023
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 18023
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -18837
$ws.Range("H100").Value = 1436.4375
$ws.Range("I100").Value = 1436.4375
$ws.Range("K100").Value = 1436.4375
$ws.Range("M100").Value = -895.4375

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 6384.6665
$ws.Range("I54").Value = 5000
$ws.Range("K54").Value = 5000
$ws.Range("M54").Value = -4480
$ws.Range("H96").Value = 1950
$ws.Range("J96").Value = 2100
$ws.Range("L96").Value = 2100
$ws.Range("N96").Value = -4846
